# Update the GESS model with the latest data:
#  - Timestamps in column A shift forward by 5 days (20/21 Feb 2026 -> 25/26 Feb 2026)
#  - The "Lookup" text in column E is rebuilt to match the new dates

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# --- Column A: shift every timestamp forward by 5 days -----------------
$rngA = $ws.Range("A2:A$lastRow")
$valsA = $rngA.Value2
$rows = $valsA.GetLength(0)
for ($i = 1; $i -le $rows; $i++) {
    $valsA[$i, 1] = $valsA[$i, 1] + 5
}
$rngA.Value2 = $valsA

# --- Column E: rewrite the "Lookup" strings for the new dates ----------
$rngE = $ws.Range("E2:E$lastRow")
$valsE = $rngE.Value2
for ($i = 1; $i -le $rows; $i++) {
    $txt = [string]$valsE[$i, 1]
    $txt = $txt.Replace("20.02.2026", "25.02.2026")
    $txt = $txt.Replace("21.02.2026", "26.02.2026")
    $valsE[$i, 1] = $txt
}
$rngE.Value2 = $valsE

Write-Output "Updated $rows data rows"
